$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.4540694355964661
$ws.Range("B1").Value = 0.4067824184894562
$ws.Range("C1").Value = 4.516948699951172
$ws.Range("D1").Value = 2.694506168365479
$ws.Range("E1").Value = 1.211387395858765
